# Weekly update: a new price record was added for "Papa" (Macroferia
# Regional de Talca) on top of the existing history table. Excel-wise this
# is simply inserting one new row right before the current row 858 (which
# pushes the whole 858:892 block down to 859:893) and then filling the new
# row with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 858; everything currently at 858:892 shifts
# down to 859:893 (values, styles and formats all move with the cells).
$ws.Rows("858").Insert()

# Populate the new row 858 with the new observation. Columns A, B, C, E,
# F, G, Q and R repeat the constant values used throughout this block.
$ws.Range("A858").Value = 5
$ws.Range("B858").Value = "Macroferia Regional de Talca"
$ws.Range("C858").Value = "Maule"
$ws.Range("D858").Value = 45147
$ws.Range("E858").Value = 7
$ws.Range("F858").Value = 100114001
$ws.Range("G858").Value = "Papa"
$ws.Range("H858").Value = "Asterix"
$ws.Range("I858").Value = "1a (cosecha)"
$ws.Range("J858").Value = 1500
$ws.Range("K858").Value = 17000
$ws.Range("L858").Value = 17000
$ws.Range("M858").Value = 17000
$ws.Range("N858").Value = "$/saco 25 kilos"
$ws.Range("O858").Value = "Región del Maule"
$ws.Range("P858").Value = 680
$ws.Range("Q858").Value = 25
$ws.Range("R858").Value = "Hortaliza"
